$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas_criteria")

$ws.Range("B2").Value = 830
$ws.Range("C2").Value = 41.9

$ws.Range("B3").Value = 1465
$ws.Range("C3").Value = 73.90000000000001

$ws.Range("B4").Value = 1229
$ws.Range("C4").Value = 62

$ws.Range("B5").Value = 1667
$ws.Range("C5").Value = 84.09999999999999
